$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from H1 into the two new
# header cells, then set their text.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows for the two new columns I (I0) and J (IF).
$data = @{
    2  = @(9, 9)
    3  = @(6, 7)
    4  = @(6, 8)
    5  = @(7, 7)
    6  = @(7, 9)
    7  = @(1, 1)
    8  = @(1, 5)
    9  = @(1, 6)
    10 = @(1, 4)
    11 = @(1, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
